# gpsSearcher / switcherr integration
# - Sheet2 ("gpsSearcher" log) gains the Helsinki/St Petersburg rows that used to
#   live on Sheet1, plus a new batch of Scandinavian / Northern-Europe / Iberian
#   places discovered by the switcher.
# - Sheet1 is "switched" to hold the latest search results only: Madrid, Bilbao
#   and Valencia.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$timeStamp = $ws2.Range("C14").Value()

# ---- Sheet2: move the old Sheet1 rows down into the log ----
$ws2.Range("A15").Value = "st petersberg"
$ws2.Range("A16").Value = "helsinki"
$ws2.Range("B15").Value = "59.9311° N, 30.3609° E"
$ws2.Range("B16").Value = "60.1699° N, 24.9384° E"
$ws2.Range("C15").Value = $timeStamp
$ws2.Range("C16").Value = $timeStamp

# ---- Sheet2: Scandinavia / Baltic coast batch ----
$ws2.Range("A17").Value = "oslo"
$ws2.Range("A18").Value = "malmo"
$ws2.Range("A19").Value = "copenhagen"
$ws2.Range("A20").Value = "odense"
$ws2.Range("A21").Value = "kiel"
$ws2.Range("A22").Value = "rostock"
$ws2.Range("B17").Value = "59.9139° N, 10.7522° E"
$ws2.Range("B18").Value = "55.6050° N, 13.0038° E"
$ws2.Range("B19").Value = "55.6761° N, 12.5683° E"
$ws2.Range("B20").Value = "55.4038° N, 10.4024° E"
$ws2.Range("B21").Value = "54.3233° N, 10.1228° E"
$ws2.Range("B22").Value = "54.0924° N, 12.0991° E"
$ws2.Range("C17:C22").Value = $timeStamp

# ---- Sheet2: Arctic Norway batch ----
$ws2.Range("A23").Value = "Svalbard"
$ws2.Range("A24").Value = "Tromso"
$ws2.Range("B23").Value = "77.8750° N, 20.9752° E"
$ws2.Range("B24").Value = "69.6492° N, 18.9553° E"
$ws2.Range("C23:C24").Value = $timeStamp

# ---- Sheet2: Germany / UK batch ----
$ws2.Range("A25").Value = "Hannover"
$ws2.Range("A26").Value = "Hamburg"
$ws2.Range("A27").Value = "London"
$ws2.Range("B25").Value = "52.3759° N, 9.7320° E"
$ws2.Range("B26").Value = "53.5511° N, 9.9937° E"
$ws2.Range("B27").Value = "51.5074° N, 0.1278° W"
$ws2.Range("C25:C27").Value = $timeStamp

# ---- Sheet2: Lisbon ----
$ws2.Range("A28").Value = "Lisbon"
$ws2.Range("B28").Value = "38.7223° N, 9.1393° W"
$ws2.Range("C28").Value = $timeStamp

# ---- Sheet2: Spain batch (also becomes the new Sheet1 content) ----
$ws2.Range("A29").Value = "Madrid"
$ws2.Range("A30").Value = "bilbao"
$ws2.Range("A31").Value = "valencia"
$ws2.Range("B29").Value = "40.4168° N, 3.7038° W"
$ws2.Range("B30").Value = "43.2630° N, 2.9350° W"
$ws2.Range("B31").Value = "39.4699° N, 0.3763° W"
$ws2.Range("C29:C31").Value = $timeStamp

# ---- Sheet1: switcherr swaps in the newest search results ----
$ws1.Range("A2").Value = "Madrid"
$ws1.Range("A3").Value = "bilbao"
$ws1.Range("A4").Value = "valencia"
$ws1.Range("B2").Value = "40.4168° N, 3.7038° W"
$ws1.Range("B3").Value = "43.2630° N, 2.9350° W"
$ws1.Range("B4").Value = "39.4699° N, 0.3763° W"

# ---- Selections: Sheet1 lands on A4, Sheet2 (still the active tab) on C22:C31 ----
$ws1.Range("A4").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("C22:C31").Select() | Out-Null
